# New Submission Synced: 2026-02-04 19:19:08
# Target sheet is "JSS 3F" (the sheet holding the "Zara Muhammad" row).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3F")

# C2's "Admission No" value was written as text ("38"); convert it to a
# real number to match the new export format.
$ws.Range("C2").Value = 38

# Append the newly synced submission as row 3.
$ws.Range("A3").Value = "2026-02-04 19:19:08"
$ws.Range("B3").Value = "ANNABEL JOEL "

# "Admission No" stays text-typed like the other rows, so force a text
# entry (leading apostrophe = Excel's "treat as text" convention) and
# then reset the cell style back to Normal so we don't leave a stray
# quote-prefixed style behind on the cell.
$ws.Range("C3").Value = "'28"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = 8
